$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the existing row 2, pushing the old row 2 data down to row 3
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).ClearFormats()

# Populate the newly inserted row 2 with the new entry (38602_lxx / 38602_mt)
$ws.Range("A2").Value = "38602_lxx"
$ws.Range("B2").Value = "38602_mt"
$ws.Range("C2").Value = "ἐξῆλθε δὲ Μωυσῆς εἰς συνάντησιν τῷ γαμβρῷ καὶ προσεκύνησεν αὐτῷ καὶ ἐφίλησεν αὐτόν, καὶ ἠσπάσαντο ἀλλήλους · καὶ εἰσήγαγεν αὐτοὺς εἰς τὴν σκηνήν."
$ws.Range("E2").Value = "(18, 7)"
$ws.Range("F2").Value = "εἰς"
$ws.Range("G2").Value = "W JY> MCH L QR>T XTNW W JCTXW W JCQ LW W JC>LW >JC L R<HW L CLWM W JB>W H >HLH"
$ws.Range("H2").Value = "H >HLH"
$ws.Range("I2").Value = "dir-he"
$ws.Range("K2").Value = "inanim"
